# Update "想去人数" (attendance count) figures on both the "展览" and
# "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3084
$ws1.Range("F6").Value = 2059
$ws1.Range("F8").Value = 145
$ws1.Range("F11").Value = 909

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3084
$ws4.Range("F6").Value = 2059
$ws4.Range("F9").Value = 145
$ws4.Range("F12").Value = 909
